$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price (D) and 1h volume change (E) figures.
# A leading apostrophe is used for numeric-looking price strings so Excel
# stores them as text (matching the original inline-string cell type)
# instead of silently converting them to numbers.

$ws.Range("D2").Value = "26.360.77"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "1.795.06"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'307.42"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").Value = "'0.4498"
$ws.Range("E7").Value = "  -1.71%  "
$ws.Range("D8").Value = "'0.3598"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "'45.83"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'0.07072"
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").Value = "'0.8844"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "'0.07747"
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "'19.38"
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").Value = "1.772.42"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "'5.289"
$ws.Range("E15").Value = "  -0.93%  "
$ws.Range("D16").Value = "'6.337"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D17").Value = "'84.84"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'0.000008499"
$ws.Range("E19").Value = "  -2.52%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "26.391.28"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").Value = "'4.983"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "'10.54"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "2.007.14"
$ws.Range("E25").Value = "  -2.00%  "
$ws.Range("D26").Value = "'1.974"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").Value = "'151.93"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'17.84"
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'2.024"
$ws.Range("E29").Value = "  +2.64%  "
$ws.Range("D30").Value = "'112.19"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").Value = "'4.870"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").Value = "'0.08685"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "'3.058"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "'2.738"
$ws.Range("E34").Value = "  +6.14%  "
$ws.Range("D35").Value = "'4.439"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").Value = "'0.7241"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "'1.105"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'1.067"
$ws.Range("E39").Value = "  -2.07%  "
$ws.Range("D40").Value = "'0.01930"
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'0.05096"
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("D42").Value = "'2.860"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'6.894"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "'0.5066"
$ws.Range("E44").Value = "  +1.81%  "
$ws.Range("D45").Value = "'0.1510"
$ws.Range("D46").Value = "'8.006"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").Value = "'1.007"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("D48").Value = "'0.4622"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").Value = "'101.29"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'9.822"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").Value = "'1.581"
$ws.Range("E51").Value = "  -2.07%  "
